$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lower-case the "MW" suffix in the header labels for the onshore
# wind/solar and storage columns (casing/wording fix).
$ws.Range("D1").Value = "apco_onshore_wind_and_solar_mw"
$ws.Range("E1").Value = "dominion_onshore_wind_and_solar_mw"
$ws.Range("F1").Value = "apco_storage_mw"
$ws.Range("G1").Value = "dominion_storage_mw"

# Move the active selection from C10 to A3.
$ws.Range("A3").Select()
